# feat: add 2022-Q4 data
#
# A new quarterly snapshot ("2022-Q4") is inserted as the newest sheet,
# right after the summary sheet "总计". The sheet that used to be the
# newest ("2022-Q3") is duplicated: the duplicate becomes the new
# "2022-Q4" sheet (with refreshed fund metrics), while the original
# keeps its data unchanged and simply becomes the (now second-newest)
# "2022-Q3" sheet. The other historical sheets ("2021-Q4", "2020-Q4")
# are unaffected aside from sliding one position to the right.
# The "总计" (summary) sheet gains a new top row for 2022-Q4 and all
# other rows shift down by one, with a brand-new row appended for the
# quarter that fell out of the 3-quarter window.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Duplicate the current newest data sheet ("2022-Q3") to create
#    the new "2022-Q4" sheet, inserted right before it.
# ---------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item("2022-Q3")
$oldQ3.Copy($oldQ3)
$newQ4 = $wb.Worksheets.Item("2022-Q3 (2)")
$newQ4.Name = "2022-Q4"

# Refresh the fund metrics on the new "2022-Q4" sheet (fund codes /
# names / fund size / G column stay the same as last quarter).
$newQ4.Range("D2").Value = "'0.16"
$newQ4.Range("E2").Value = "'93.15"
$newQ4.Range("F2").Value = "'1.21"
$newQ4.Range("H2").Value = 10

$newQ4.Range("E3").Value = "'93.15"
$newQ4.Range("F3").Value = "'1.21"
$newQ4.Range("H3").Value = 10

# ---------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new top data row
#    for 2022-Q4 and append a new row for 2020-Q4 at the bottom.
# ---------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

$zj.Rows.Item(2).Insert()

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 0

# The insert leaves row 2 with borrowed formatting; re-apply the
# standard index-column style (copied from the row below, which kept
# the original formatting of what used to be row 2).
$zj.Range("A2:D2").ClearFormats()
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

# Renumber the index column for the rows that shifted down.
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
$zj.Range("A5").Value = 3

# Keep the summary sheet as the active tab, same as before the edit.
$zj.Activate() | Out-Null
